$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 383, shifting existing rows 383:435 down to 384:436
$ws.Rows.Item(383).Insert()

# Populate the newly inserted row 383 with the new data record
$ws.Range("A383").Value = 4
$ws.Range("B383").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C383").Value = "Los Lagos"
$ws.Range("D383").Value = 45131
$ws.Range("E383").Value = 10
$ws.Range("F383").Value = 100112032
$ws.Range("G383").Value = "Zapallo italiano"
$ws.Range("H383").Value = "Sin especificar"
$ws.Range("I383").Value = "Primera"
$ws.Range("J383").Value = 250
$ws.Range("K383").Value = 20000
$ws.Range("L383").Value = 20000
$ws.Range("M383").Value = 20000
$ws.Range("N383").Value = "$/caja 50 unidades"
$ws.Range("O383").Value = "Región de Arica y Parinacota"
$ws.Range("P383").Value = 400
$ws.Range("Q383").Value = 50
$ws.Range("R383").Value = "Hortaliza"

Write-Host "Done"
